# Update column G ("K") values on the active sheet, rows 2-35, to reflect
# the regenerated save_data (K instead of Strike#, recalculated std/mean,
# s_vals written out).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 5
    3  = 6
    4  = 3
    5  = 1
    6  = 5
    7  = 4
    8  = 6
    9  = 11
    10 = 2
    11 = 1
    12 = 6
    13 = 4
    14 = 4
    15 = 8
    16 = 4
    17 = 2
    18 = 4
    19 = 5
    20 = 4
    21 = 6
    22 = 3
    23 = 5
    24 = 2
    25 = 4
    26 = 4
    27 = 4
    28 = 4
    29 = 4
    30 = 2
    31 = 4
    32 = 5
    33 = 2
    34 = 6
    35 = 1
}

foreach ($row in $newValues.Keys | Sort-Object) {
    $ws.Range("G$row").Value = $newValues[$row]
}
